$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row before row 145, shifting rows 145:234 down to 146:235
$ws.Rows.Item(145).Insert()

# Populate the new row 145 with values
$ws.Cells.Item(145, 1).Value = 3
$ws.Cells.Item(145, 2).Value = "Femacal de La Calera"
$ws.Cells.Item(145, 3).Value = "Coquimbo"
$ws.Cells.Item(145, 4).Value = 44603
$ws.Cells.Item(145, 5).Value = 5
$ws.Cells.Item(145, 6).Value = 100112001
$ws.Cells.Item(145, 7).Value = "Berenjena"
$ws.Cells.Item(145, 8).Value = "Sin especificar"
$ws.Cells.Item(145, 9).Value = "Primera"
$ws.Cells.Item(145, 10).Value = 85
$ws.Cells.Item(145, 11).Value = 9000
$ws.Cells.Item(145, 12).Value = 9500
$ws.Cells.Item(145, 13).Value = 9235
$ws.Cells.Item(145, 14).Value = "$/caja 60 unidades"
$ws.Cells.Item(145, 15).Value = "Región Metropolitana"
$ws.Cells.Item(145, 16).Value = 154
$ws.Cells.Item(145, 17).Value = 60
$ws.Cells.Item(145, 18).Value = "Hortaliza"
